$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The time-entry rows got re-sorted into a new order (the download query
# now orders results differently). Row 3 (idtimeentry 2) is unaffected and
# is left completely untouched. For the other rows, column A holds the
# numeric-looking idtimeentry id as text, so force a text number format
# before writing the value back -- otherwise "5" etc. would be silently
# coerced to a real number by Excel.
$rows = @(2, 4, 5, 6, 7)
foreach ($r in $rows) {
    $ws.Cells.Item($r, 1).NumberFormat = "@"
}

$data = @{
    2 = @("5", "2024-05-08 00:00:00", "4:40:00",  "5:41:00",  "2024-05-22 01:37:38", "2024-05-22 01:37:38", "johndoe")
    4 = @("6", "2024-05-21 00:00:00", "15:57:00", "18:00:00", "2024-05-22 01:55:36", "2024-05-22 01:55:36", "avisemah")
    5 = @("1", "2024-05-22 00:00:00", "4:07:00",  "7:10:00",  "2024-05-22 01:04:37", "2024-05-22 01:04:37", "johndoe")
    6 = @("3", "2024-05-27 00:00:00", "4:29:00",  "7:32:00",  "2024-05-22 01:26:23", "2024-05-22 01:26:23", "johndoe")
    7 = @("4", "2024-05-28 00:00:00", "3:36:00",  "5:38:00",  "2024-05-22 01:35:20", "2024-05-22 01:35:20", "johndoe")
}

foreach ($row in $rows) {
    $rec = $data[$row]
    $ws.Cells.Item($row, 1).Value = $rec[0]
    $ws.Cells.Item($row, 2).Value = $rec[1]
    $ws.Cells.Item($row, 3).Value = $rec[2]
    $ws.Cells.Item($row, 4).Value = $rec[3]
    $ws.Cells.Item($row, 5).Value = $rec[4]
    $ws.Cells.Item($row, 6).Value = $rec[5]
    $ws.Cells.Item($row, 7).Value = $rec[6]
}
